# Iraq League - Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
#
# For a handful of match rows, the row that should have held one fixture's
# data actually held the other's (and vice versa). The fix swaps the full
# row contents (every column except the leading id column A, which must
# stay attached to its own row number) between each of the following
# row pairs: (17,18) (69,70) (89,90) (135,136) (219,220) (223,224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(17, 18),
    @(69, 70),
    @(89, 90),
    @(135, 136),
    @(219, 220),
    @(223, 224)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Columns B (2) through AD (30) hold the swappable data; column A holds
    # the row's own id and must not move.
    $range1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 30))
    $range2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 30))

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
